# Regenerate the s_vals data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) and the computed sum
# column G for each data row (rows 2-10) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    @{ Row = 3;  B = 1.455362044514542;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 3.754798637575387 }
    @{ Row = 4;  B = 0.2917716402565462; C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 3.194529526351451 }
    @{ Row = 5;  B = 0.1190320826869504; C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 3.021789968781855 }
    @{ Row = 6;  B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    @{ Row = 7;  B = 0.6606524410359556; C = 0.306821227259698;  D = 0.7527432677738641; E = 0.4942365360607697; G = 2.214453472130288 }
    @{ Row = 8;  B = 1.455362044514542;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 3.754798637575387 }
    @{ Row = 9;  B = 0.1190320826869504; C = 0.04071648406533734; D = 0.7527432677738641; E = 0.4942365360607697; G = 1.406728370586922 }
    @{ Row = 10; B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 7).Value = $entry.G
}
